$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 5 new rows at the top of the block (rows 1073-1077), shifting
# the existing rows 1073-1137 down to 1078-1142.
$ws.Rows("1073:1077").Insert()

# Row 1073
$ws.Range("A1073").Value = 3
$ws.Range("B1073").Value = 'Femacal de La Calera'
$ws.Range("C1073").Value = 'Coquimbo'
$ws.Range("D1073").Value = 44516
$ws.Range("E1073").Value = 5
$ws.Range("F1073").Value = 100112033
$ws.Range("G1073").Value = 'Lechuga'
$ws.Range("H1073").Value = 'Conconina(o)'
$ws.Range("I1073").Value = 'Primera'
$ws.Range("J1073").Value = 130
$ws.Range("K1073").Value = 4000
$ws.Range("L1073").Value = 4500
$ws.Range("M1073").Value = 4269
$ws.Range("N1073").Value = '$/caja 10 unidades'
$ws.Range("O1073").Value = 'Provincia de Quillota'
$ws.Range("P1073").Value = 427
$ws.Range("Q1073").Value = 10
$ws.Range("R1073").Value = 'Hortaliza'

# Row 1074
$ws.Range("A1074").Value = 3
$ws.Range("B1074").Value = 'Femacal de La Calera'
$ws.Range("C1074").Value = 'Coquimbo'
$ws.Range("D1074").Value = 44516
$ws.Range("E1074").Value = 5
$ws.Range("F1074").Value = 100112033
$ws.Range("G1074").Value = 'Lechuga'
$ws.Range("H1074").Value = 'Escarola'
$ws.Range("I1074").Value = 'Primera'
$ws.Range("J1074").Value = 125
$ws.Range("K1074").Value = 5500
$ws.Range("L1074").Value = 5800
$ws.Range("M1074").Value = 5656
$ws.Range("N1074").Value = '$/caja 15 unidades'
$ws.Range("O1074").Value = 'Provincia de Quillota'
$ws.Range("P1074").Value = 377
$ws.Range("Q1074").Value = 15
$ws.Range("R1074").Value = 'Hortaliza'

# Row 1075
$ws.Range("A1075").Value = 3
$ws.Range("B1075").Value = 'Femacal de La Calera'
$ws.Range("C1075").Value = 'Coquimbo'
$ws.Range("D1075").Value = 44516
$ws.Range("E1075").Value = 5
$ws.Range("F1075").Value = 100112033
$ws.Range("G1075").Value = 'Lechuga'
$ws.Range("H1075").Value = 'Francesa morada'
$ws.Range("I1075").Value = 'Primera'
$ws.Range("J1075").Value = 60
$ws.Range("K1075").Value = 5000
$ws.Range("L1075").Value = 5000
$ws.Range("M1075").Value = 5000
$ws.Range("N1075").Value = '$/caja 18 unidades'
$ws.Range("O1075").Value = 'Provincia de Quillota'
$ws.Range("P1075").Value = 278
$ws.Range("Q1075").Value = 18
$ws.Range("R1075").Value = 'Hortaliza'

# Row 1076
$ws.Range("A1076").Value = 3
$ws.Range("B1076").Value = 'Femacal de La Calera'
$ws.Range("C1076").Value = 'Coquimbo'
$ws.Range("D1076").Value = 44516
$ws.Range("E1076").Value = 5
$ws.Range("F1076").Value = 100112033
$ws.Range("G1076").Value = 'Lechuga'
$ws.Range("H1076").Value = 'Marina'
$ws.Range("I1076").Value = 'Primera'
$ws.Range("J1076").Value = 60
$ws.Range("K1076").Value = 5000
$ws.Range("L1076").Value = 5000
$ws.Range("M1076").Value = 5000
$ws.Range("N1076").Value = '$/caja 18 unidades'
$ws.Range("O1076").Value = 'Provincia de Quillota'
$ws.Range("P1076").Value = 278
$ws.Range("Q1076").Value = 18
$ws.Range("R1076").Value = 'Hortaliza'

# Row 1077
$ws.Range("A1077").Value = 3
$ws.Range("B1077").Value = 'Femacal de La Calera'
$ws.Range("C1077").Value = 'Coquimbo'
$ws.Range("D1077").Value = 44516
$ws.Range("E1077").Value = 5
$ws.Range("F1077").Value = 100112033
$ws.Range("G1077").Value = 'Lechuga'
$ws.Range("H1077").Value = 'Milanesa'
$ws.Range("I1077").Value = 'Primera'
$ws.Range("J1077").Value = 130
$ws.Range("K1077").Value = 4500
$ws.Range("L1077").Value = 4800
$ws.Range("M1077").Value = 4650
$ws.Range("N1077").Value = '$/caja 20 unidades'
$ws.Range("O1077").Value = 'Provincia de Quillota'
$ws.Range("P1077").Value = 232
$ws.Range("Q1077").Value = 20
$ws.Range("R1077").Value = 'Hortaliza'
